# Applies the changes described in the commit:
# "Added more hand-verified sequence sheets"
#
# - Renames the worksheet from "Sheet1" to "Sheet"
# - Fills in the newly added "A1" label in column C for rows 2-4
# - Updates the view: zoom level and the active/selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Sheet"

# Populate the new column C values (label referencing cell A1)
$ws.Range("C2").Value = "A1"
$ws.Range("C3").Value = "A1"
$ws.Range("C4").Value = "A1"

# Update the window zoom level
$win = $excel.ActiveWindow
$win.Zoom = 215

# Update the selected/active cell shown in the sheet view
[void]$ws.Range("B6").Select()
